# Apply the edit described by the commit:
#  - Modify slide "Design Strategy - Connection Pooling" (SlideID 533):
#      * move/resize the Content Placeholder 2 shape
#      * change its text run "that maintains several threads. " -> "that maintains several connections to support several threads."
#      * move the picture shape
#  - Delete the two exercise solution slides ("EXERCISE - R Solution", "EXERCISE - Java Solution")
#  - Update the notes-page cached slide-number field on the slide that shifts from #19 to #17
#  - Refresh the cached "last saved" date fields (10/1/2020 -> 2/22/2021) on the
#    handout master, notes master, slide master and the "Title and Content" layout

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Edit the "Design Strategy - Connection Pooling" slide (currently slide 22,
#    SlideID 533) BEFORE deleting any slides so indices stay simple.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 533) {
        $target = $p.Slides.Item($i)
        break
    }
}

$content = $target.Shapes.Item(2)
$pic = $target.Shapes.Item(3)

# Reposition / resize the content placeholder
$content.Left = 539494 / 12700.0
$content.Top = 1435608 / 12700.0
$content.Width = 611.7992
$content.Height = 5277708 / 12700.0

# Update the text run in-place, preserving its run formatting
$tr = $content.TextFrame.TextRange
$full = $tr.Text
$oldRun = "that maintains several threads. "
$newRun = "that maintains several connections to support several threads."
$idx = $full.IndexOf($oldRun)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldRun.Length)
    $sub.Text = $newRun
}

# Reposition the picture
$pic.Left = 8309343 / 12700.0
$pic.Top = 1435608 / 12700.0

# ---------------------------------------------------------------------------
# 2. Update the cached notes-page slide-number field for the slide that will
#    shift from position 19 to position 17 once the two slides below are
#    removed (SlideID 522, currently at index 19).
# ---------------------------------------------------------------------------
$notesTarget = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 522) {
        $notesTarget = $p.Slides.Item($i)
        break
    }
}
if ($notesTarget -ne $null -and $notesTarget.HasNotesPage) {
    $np = $notesTarget.NotesPage
    for ($j = 1; $j -le $np.Shapes.Count; $j++) {
        $sh = $np.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "19") {
            $sh.TextFrame.TextRange.Text = "17"
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Delete the two "EXERCISE" solution slides (SlideID 529 "R Solution" and
#    SlideID 530 "Java Solution").
# ---------------------------------------------------------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $sid = $p.Slides.Item($i).SlideID
    if ($sid -eq 529 -or $sid -eq 530) {
        $p.Slides.Item($i).Delete()
    }
}

# ---------------------------------------------------------------------------
# 4. Refresh the cached date fields across the masters/layout.
# ---------------------------------------------------------------------------
$newDate = "2/22/2021"

function Update-DateField($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $s = $shapes.Item($k)
        if ($s.HasTextFrame -and $s.TextFrame.HasText -and $s.TextFrame.TextRange.Text -eq "10/1/2020") {
            $s.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DateField $p.HandoutMaster.Shapes
Update-DateField $p.NotesMaster.Shapes
Update-DateField $p.SlideMaster.Shapes
Update-DateField $p.SlideMaster.CustomLayouts.Item(2).Shapes
